$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "success" column (D): header in D1 matching the style used
# by the other header cells (B1/C1), and a 0/1 flag per data row.

$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats - reuse existing style, no new style record
$ws.Range("D1").Value = "success"

$successValues = @("0", "0", "0", "0", "0", "1", "0", "0", "0", "0", "0", "0")

for ($i = 0; $i -lt $successValues.Length; $i++) {
    $r = $i + 2
    $cell = $ws.Cells.Item($r, 4)
    # Route the literal text through a formula + paste-as-values so Excel
    # stores it as a genuine text cell (t="s") without forcing a new
    # number-format / style record onto the cell.
    $cell.Formula = "=""" + $successValues[$i] + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)         # xlPasteValues
}

$excel.CutCopyMode = 0
